$p = $ppt.ActivePresentation

# 1) Refresh the cached "datetimeFigureOut" date field text wherever it
#    appears (Date Placeholder shapes on slides/layouts/masters) from
#    "04-08-2024" to "20-10-2024". Harmless no-op on slides that don't
#    have such a field.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "04-08-2024") {
                $tr.Text = "20-10-2024"
            }
        }
    }
}

# 2) On the title/index slide, the "UNIT-1" arrow shape currently holds
#    two paragraphs: "UNIT-1" and "Introduction to Computers". Drop the
#    second paragraph, leaving just "UNIT-1" (matching the other
#    "UNIT-n" arrow shapes, which only ever carried one line).
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $full = $shp.TextFrame.TextRange.Text
            if ($full.StartsWith("UNIT-1") -and $full.Contains("Introduction to Computers")) {
                $shp.TextFrame.TextRange.Text = "UNIT-1"
            }
        }
    }
}
